# Edit: split the run
#   "Source of our data is AACT database containing all the "
# into the following sequence of runs (text content only -- the fine
# grained rPr markers such as lang="en-US"/err="1" that PowerPoint's UI
# leaves behind when text is retyped/pasted are not reachable through the
# COM object model and are left to whatever the host assigns by default):
#   "Source of our data is "
#   " the Aggregate Analysis of "
#   "ClinicalTrials.gov"
#   " ("
#   "AACT"
#   ")"
#   " "
#   " "
#   "database "
#   "containing all the "
# while leaving the following two runs ("ClinicalTrials.gov" / " records")
# untouched.

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
$targetParaIndex = 0
$needle = "AACT database"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
  $sl = $p.Slides.Item($si)
  for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
    $shp = $sl.Shapes.Item($shi)
    if ($shp.HasTextFrame) {
      $tf = $shp.TextFrame
      if ($tf.HasText) {
        $tr = $tf.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
          $para = $tr.Paragraphs($pi, 1)
          if ($para.Text -like "*$needle*") {
            $targetSlide = $sl
            $targetShape = $shp
            $targetParaIndex = $pi
          }
        }
      }
    }
  }
}

if ($targetShape -eq $null) {
  throw "Could not find the paragraph containing '$needle'"
}

$tr = $targetShape.TextFrame.TextRange
$para = $tr.Paragraphs($targetParaIndex, 1)

# New chunks that together replace the original first run's text
# ("Source of our data is AACT database containing all the ").
$chunks = @(
  "Source of our data is ",
  " the Aggregate Analysis of ",
  "ClinicalTrials.gov",
  " (",
  "AACT",
  ")",
  " ",
  " ",
  "database ",
  "containing all the "
)

# The paragraph-level InsertBefore() always prepends right at the start of
# the paragraph (regardless of which sub-range it is invoked on) and the
# newly created run inherits the rPr of whichever run currently sits at
# the front of the paragraph. Since that front run is (and remains,
# throughout this loop) the original first run -- the one we are
# splitting -- every newly inserted chunk correctly inherits that run's
# plain `dirty="0"` formatting instead of bleeding in formatting from
# neighbouring runs. Insert the chunks in reverse order so they end up in
# the right left-to-right order once all insertions are done.
for ($i = $chunks.Length - 1; $i -ge 0; $i--) {
  $para.InsertBefore($chunks[$i]) | Out-Null
}

# After the loop the paragraph looks like:
#   [chunk1]...[chunk10][original first run, unchanged]["ClinicalTrials.gov"][" records"]
# The original first run content is now fully redundant (it has been
# replaced piece-by-piece by chunk1..chunk10), so clear it out, which
# removes that now-empty run entirely.
$leftoverIndex = $chunks.Length + 1
$leftover = $para.Runs($leftoverIndex, 1)
$leftover.Text = ""
